$wb = $excel.ActiveWorkbook

# ----- Sheet 1: LP1912 -----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range('A2').Value = 'Última actualización: 10:04:30'
$ws1.Range('A3').Value = 'Total filas: 100'

$ws1.Cells.Item(66, 1).Value = '08:11:18'
$ws1.Cells.Item(66, 2).Value = '09:28'
$ws1.Cells.Item(66, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(66, 4).Value = 77
$ws1.Cells.Item(66, 5).Value = 'LP1912'

$ws1.Cells.Item(67, 1).Value = '08:28:52'
$ws1.Cells.Item(67, 2).Value = '09:28'
$ws1.Cells.Item(67, 3).Value = '10_OLMOS'
$ws1.Cells.Item(67, 4).Value = 60
$ws1.Cells.Item(67, 5).Value = 'LP1912'

$ws1.Cells.Item(84, 1).Value = '10:04:30'
$ws1.Cells.Item(84, 2).Value = '10:22'
$ws1.Cells.Item(84, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(84, 4).Value = 18
$ws1.Cells.Item(84, 5).Value = 'LP1912'

$ws1.Cells.Item(85, 1).Value = '09:22:34'
$ws1.Cells.Item(85, 2).Value = '10:25'
$ws1.Cells.Item(85, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(85, 4).Value = 63
$ws1.Cells.Item(85, 5).Value = 'LP1912'

$ws1.Cells.Item(86, 1).Value = '08:38:24'
$ws1.Cells.Item(86, 2).Value = '10:29'
$ws1.Cells.Item(86, 3).Value = '15_ABASTO'
$ws1.Cells.Item(86, 4).Value = 111
$ws1.Cells.Item(86, 5).Value = 'LP1912'

$ws1.Cells.Item(87, 1).Value = '10:04:30'
$ws1.Cells.Item(87, 2).Value = '10:29'
$ws1.Cells.Item(87, 3).Value = '14_ABASTO'
$ws1.Cells.Item(87, 4).Value = 25
$ws1.Cells.Item(87, 5).Value = 'LP1912'

$ws1.Cells.Item(88, 1).Value = '08:45:31'
$ws1.Cells.Item(88, 2).Value = '10:44'
$ws1.Cells.Item(88, 3).Value = '11X44_ETCHEVERRY'
$ws1.Cells.Item(88, 4).Value = 119
$ws1.Cells.Item(88, 5).Value = 'LP1912'

$ws1.Cells.Item(89, 1).Value = '08:52:40'
$ws1.Cells.Item(89, 2).Value = '10:46'
$ws1.Cells.Item(89, 3).Value = '15_P INDUSTRIAL'
$ws1.Cells.Item(89, 4).Value = 114
$ws1.Cells.Item(89, 5).Value = 'LP1912'

$ws1.Cells.Item(90, 1).Value = '09:22:34'
$ws1.Cells.Item(90, 2).Value = '10:53'
$ws1.Cells.Item(90, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(90, 4).Value = 91
$ws1.Cells.Item(90, 5).Value = 'LP1912'

$ws1.Cells.Item(91, 1).Value = '10:04:30'
$ws1.Cells.Item(91, 2).Value = '10:56'
$ws1.Cells.Item(91, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(91, 4).Value = 52
$ws1.Cells.Item(91, 5).Value = 'LP1912'

$ws1.Cells.Item(92, 1).Value = '09:22:34'
$ws1.Cells.Item(92, 2).Value = '10:57'
$ws1.Cells.Item(92, 3).Value = '10_OLMOS'
$ws1.Cells.Item(92, 4).Value = 95
$ws1.Cells.Item(92, 5).Value = 'LP1912'

$ws1.Cells.Item(93, 1).Value = '10:04:30'
$ws1.Cells.Item(93, 2).Value = '10:59'
$ws1.Cells.Item(93, 3).Value = '10_OLMOS'
$ws1.Cells.Item(93, 4).Value = 55
$ws1.Cells.Item(93, 5).Value = 'LP1912'

$ws1.Cells.Item(94, 1).Value = '09:22:34'
$ws1.Cells.Item(94, 2).Value = '11:01'
$ws1.Cells.Item(94, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(94, 4).Value = 99
$ws1.Cells.Item(94, 5).Value = 'LP1912'

$ws1.Cells.Item(95, 1).Value = '10:04:30'
$ws1.Cells.Item(95, 2).Value = '11:03'
$ws1.Cells.Item(95, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(95, 4).Value = 59
$ws1.Cells.Item(95, 5).Value = 'LP1912'

$ws1.Cells.Item(96, 1).Value = '09:22:34'
$ws1.Cells.Item(96, 2).Value = '11:10'
$ws1.Cells.Item(96, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(96, 4).Value = 108
$ws1.Cells.Item(96, 5).Value = 'LP1912'

$ws1.Cells.Item(97, 1).Value = '09:22:34'
$ws1.Cells.Item(97, 2).Value = '11:14'
$ws1.Cells.Item(97, 3).Value = '14_ABASTO'
$ws1.Cells.Item(97, 4).Value = 112
$ws1.Cells.Item(97, 5).Value = 'LP1912'

$ws1.Cells.Item(98, 1).Value = '09:22:34'
$ws1.Cells.Item(98, 2).Value = '11:15'
$ws1.Cells.Item(98, 3).Value = '15X38_ABASTO'
$ws1.Cells.Item(98, 4).Value = 113
$ws1.Cells.Item(98, 5).Value = 'LP1912'

$ws1.Cells.Item(99, 1).Value = '10:04:30'
$ws1.Cells.Item(99, 2).Value = '11:29'
$ws1.Cells.Item(99, 3).Value = '10_OLMOS'
$ws1.Cells.Item(99, 4).Value = 85
$ws1.Cells.Item(99, 5).Value = 'LP1912'

$ws1.Cells.Item(100, 1).Value = '10:04:30'
$ws1.Cells.Item(100, 2).Value = '11:29'
$ws1.Cells.Item(100, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(100, 4).Value = 85
$ws1.Cells.Item(100, 5).Value = 'LP1912'

$ws1.Cells.Item(101, 1).Value = '10:04:30'
$ws1.Cells.Item(101, 2).Value = '11:31'
$ws1.Cells.Item(101, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(101, 4).Value = 87
$ws1.Cells.Item(101, 5).Value = 'LP1912'

$ws1.Cells.Item(102, 1).Value = '10:04:30'
$ws1.Cells.Item(102, 2).Value = '11:41'
$ws1.Cells.Item(102, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(102, 4).Value = 97
$ws1.Cells.Item(102, 5).Value = 'LP1912'

$ws1.Cells.Item(103, 1).Value = '10:04:30'
$ws1.Cells.Item(103, 2).Value = '11:45'
$ws1.Cells.Item(103, 3).Value = '15X38_ABASTO'
$ws1.Cells.Item(103, 4).Value = 101
$ws1.Cells.Item(103, 5).Value = 'LP1912'

$ws1.Cells.Item(104, 1).Value = '10:04:30'
$ws1.Cells.Item(104, 2).Value = '11:53'
$ws1.Cells.Item(104, 3).Value = '225_GOMEZ'
$ws1.Cells.Item(104, 4).Value = 109
$ws1.Cells.Item(104, 5).Value = 'LP1912'

$ws1.Cells.Item(105, 1).Value = '10:04:30'
$ws1.Cells.Item(105, 2).Value = '11:58'
$ws1.Cells.Item(105, 3).Value = '17_ROMERO'
$ws1.Cells.Item(105, 4).Value = 114
$ws1.Cells.Item(105, 5).Value = 'LP1912'

# ----- Sheet 2: LP1912-215 -----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range('A2').Value = 'Última actualización: 10:04:30'
$ws2.Range('A3').Value = 'Total filas: 20'

$ws2.Cells.Item(24, 1).Value = '10:04:30'
$ws2.Cells.Item(24, 2).Value = '11:31'
$ws2.Cells.Item(24, 3).Value = '215C_EL PATO'
$ws2.Cells.Item(24, 4).Value = 87
$ws2.Cells.Item(24, 5).Value = 'LP1912'

$ws2.Cells.Item(25, 1).Value = '10:04:30'
$ws2.Cells.Item(25, 2).Value = '11:41'
$ws2.Cells.Item(25, 3).Value = '215B_EL PATO'
$ws2.Cells.Item(25, 4).Value = 97
$ws2.Cells.Item(25, 5).Value = 'LP1912'

# ----- Sheet 3: 6203-6173 -----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range('A2').Value = 'Última actualización: 10:04:30'
$ws3.Range('A3').Value = 'Total filas: 18'

$ws3.Cells.Item(23, 1).Value = '10:04:30'
$ws3.Cells.Item(23, 2).Value = '11:26'
$ws3.Cells.Item(23, 3).Value = '215C_LA PLATA'
$ws3.Cells.Item(23, 4).Value = 82
$ws3.Cells.Item(23, 5).Value = 'L6203'

Write-Host "done"